# Updated symbol list (Price / Volume(1h) refresh) - GitHub Actions run.
# Each Price/Volume cell is stored as text (e.g. "292.36", "-6.55%"), so
# values are written with a leading apostrophe to force text entry and
# keep Excel from reinterpreting them as numbers/percentages.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'292.36"
$ws.Range("E2").Value = "'-6.55%"
$ws.Range("D3").Value = "'40.37"
$ws.Range("E3").Value = "'0.86%"
$ws.Range("D4").Value = "'5.010"
$ws.Range("E4").Value = "'-3.50%"
$ws.Range("D5").Value = "'0.07317"
$ws.Range("E5").Value = "'-3.47%"
$ws.Range("D6").Value = "'1.521"
$ws.Range("E6").Value = "'-9.63%"
$ws.Range("D7").Value = "'0.9254"
$ws.Range("E7").Value = "'-0.02%"
$ws.Range("D9").Value = "'0.1206"
$ws.Range("E9").Value = "'0.63%"
$ws.Range("D10").Value = "'0.1740"
$ws.Range("E10").Value = "'-4.14%"
$ws.Range("D11").Value = "'0.04313"
$ws.Range("E11").Value = "'3.56%"
$ws.Range("D12").Value = "'0.08597"
$ws.Range("E12").Value = "'-4.73%"
$ws.Range("E13").Value = "'0.26%"
$ws.Range("D14").Value = "'0.001265"
$ws.Range("E14").Value = "'-1.24%"
$ws.Range("D15").Value = "'0.005930"
$ws.Range("E15").Value = "'1.54%"
$ws.Range("E16").Value = "'-0.36%"
$ws.Range("D17").Value = "'4.294"
$ws.Range("E17").Value = "'-1.29%"
$ws.Range("E18").Value = "'-1.99%"
$ws.Range("D19").Value = "'7.737"
$ws.Range("E19").Value = "'1.54%"
$ws.Range("E20").Value = "'2.93%"
$ws.Range("D21").Value = "'0.2793"
$ws.Range("E21").Value = "'-0.64%"
$ws.Range("D22").Value = "'0.03927"
$ws.Range("E22").Value = "'-2.02%"
$ws.Range("D23").Value = "'0.001261"
$ws.Range("E23").Value = "'-0.72%"
$ws.Range("D24").Value = "'0.003781"
$ws.Range("E24").Value = "'-4.93%"
$ws.Range("E25").Value = "'0.80%"
$ws.Range("D38").Value = "'0.02296"
$ws.Range("E38").Value = "'-5.31%"
$ws.Range("D39").Value = "'0.04974"
$ws.Range("E39").Value = "'-3.53%"
$ws.Range("D40").Value = "'0.005414"
$ws.Range("E40").Value = "'86.90%"
$ws.Range("D41").Value = "'0.007693"
$ws.Range("E41").Value = "'-0.40%"
$ws.Range("D42").Value = "'0.1285"
$ws.Range("E42").Value = "'-1.06%"
$ws.Range("D43").Value = "'0.007326"
$ws.Range("E43").Value = "'-3.92%"
$ws.Range("D44").Value = "'0.007919"
$ws.Range("E44").Value = "'-3.61%"
$ws.Range("D45").Value = "'0.3181"
$ws.Range("E45").Value = "'2.36%"
$ws.Range("D46").Value = "'0.00006322"
$ws.Range("E46").Value = "'-3.98%"
$ws.Range("E47").Value = "'0.02%"
$ws.Range("D48").Value = "'0.02042"
$ws.Range("E48").Value = "'-92.34%"
$ws.Range("E49").Value = "'0.02%"
$ws.Range("E50").Value = "'0.02%"
